$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column F (Resultado Conseguido) wording in rows 3-10 ---
$ws.Range("F3").Value = "Te mandar a sus correspondidas secciones."

$ws.Range("F4").Value = "Se muestra:`r`n1.Background blanco a gris`r`n2.Letras de negras a blancas`r`n3.Al presionar se muestra el registro al Login"

$ws.Range("F5").Value = "Se muestra:`r`n1. Se resalta a un tono mas claro.`r`n2.Poder escribir un texto en el campo."

$ws.Range("F6").Value = "Se muestra:`r`n1.Se resalta a un tono mas claro.`r`n2.Cambia el texto a asteriscos."

$ws.Range("F7").Value = "Se muestra:`r`n1.Se resalta a un tono mas claro.`r`n2.Cambia el texto a asteriscos."

$ws.Range("F8").Value = "Se muestra:`r`nCuadro de contraseña para ingresar la contraseña"

$ws.Range("F9").Value = "Se muestra:`r`nNoticias de ultima hora"

$ws.Range("F10").Value = "`r`n1.Se resaltan las imágenes al pasar el cursor arriba."

# --- Adjust row heights that changed as a result of the re-wrapped text ---
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(7).RowHeight = 80.25
$ws.Rows.Item(10).RowHeight = 60

# --- Add the new test case row (PPDN-010) ---
$ws.Range("A11").Value = "PPDN-010"
$ws.Range("B11").Value = "Pagina principal"
$ws.Range("C11").Value = "Validar hipervinculos de las imágenes"
$ws.Range("D11").Value = "Acceso a la pagina principal"
$ws.Range("E11").Value = "Entrar a la pagina principal."
$ws.Range("F11").Value = "1.Te manda a sus correspondidas paginas"
$ws.Range("G11").Value = "Esperado"
$ws.Range("H11").Value = "Correcto"

$ws.Range("C11:H11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 45

# --- Update the selection shown when the workbook is reopened ---
$ws.Range("I11").Select() | Out-Null
